# Fiksa bruleveranse, skal ha lengde bru, ikke lengde hoydebegrensning

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Bru høydebegrensning under 4m" ---
$ws1 = $wb.Worksheets.Item(1)

# Header: "Lengde (m)" -> "Lengde"
$ws1.Range("C1").Value = "Lengde"

# Replace lengde values with correct bridge lengths
$ws1.Range("C2").Value = 118.8
$ws1.Range("C3").Value = 412.34
$ws1.Range("C4").Value = 81
$ws1.Range("C5").Value = 261.25
$ws1.Range("C6").Value = 159.82
$ws1.Range("C7").Value = 203
$ws1.Range("C8").Value = 136.18

# --- Sheet 2: "Metadata" ---
$ws2 = $wb.Worksheets.Item(2)

# Update overlapp value to note the egenskapfilter used
$ws2.Range("B5").Value = "60(1263=7304)"

# Add new metadata rows describing the property filter and the
# overlap found while searching for the height restriction
$ws2.Range("A8").Value = "egenskapfilter_bru"
$ws2.Range("B8").Value = "1263=7304"

$ws2.Range("A9").Value = "overlapp fra søk etter høydebegrensning"
$ws2.Range("B9").Value = "0.80121169-0.80129544@1125844,0.03534597-0.06032389@121713,0.97659909-0.97806091@22107,0.47770288-0.53110961@1175773,0.97700648-0.97727196@181212,0.93980412-0.94078725@22110,0.73952864-0.88387233@1060365,0.95303187-0.9588474@22110,0.87978719-0.88356741@384326,0.01180717-0.01196701@72561,0.48914175-0.88487916@1060530,0.60680381-0.60690384@805106,0.91324653-0.97386356@121475,0.96427989-0.97700648@181212,0.08871395-0.78240909@705367,0.94078725-0.95303187@22110,0.61378309-0.9735731@72630,0.21997673-0.22009309@1126325,0.98037222-0.99155114@22110,0.00109324-0.01185904@705136"
